$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "301.40"
Set-TextValue "E2" "-4.51%"
Set-TextValue "D3" "35.20"
Set-TextValue "E3" "-2.41%"
Set-TextValue "D4" "5.067"
Set-TextValue "E4" "-1.73%"
Set-TextValue "D5" "0.07966"
Set-TextValue "E5" "-2.39%"
Set-TextValue "D6" "1.916"
Set-TextValue "E6" "-10.07%"
Set-TextValue "D7" "7.733"
Set-TextValue "E7" "-3.79%"
Set-TextValue "D8" "2.911"
Set-TextValue "E8" "3.94%"
Set-TextValue "D9" "0.9230"
Set-TextValue "E10" "23.91%"
Set-TextValue "D11" "0.1853"
Set-TextValue "E11" "-1.33%"
Set-TextValue "D12" "0.09700"
Set-TextValue "E12" "5.22%"
Set-TextValue "D13" "0.03628"
Set-TextValue "E13" "0.32%"
Set-TextValue "D14" "0.09852"
Set-TextValue "E14" "-0.78%"
Set-TextValue "D15" "0.001388"
Set-TextValue "E15" "-3.29%"
Set-TextValue "D16" "0.005825"
Set-TextValue "E16" "2.31%"
Set-TextValue "D17" "3.506"
Set-TextValue "E17" "0.86%"
Set-TextValue "D18" "4.067"
Set-TextValue "E18" "-2.07%"
Set-TextValue "D19" "0.3400"
Set-TextValue "E19" "0.85%"
Set-TextValue "E20" "-1.91%"
Set-TextValue "D21" "5.062"
Set-TextValue "E21" "-2.35%"
Set-TextValue "D22" "0.2403"
Set-TextValue "E22" "9.56%"
Set-TextValue "D23" "0.04532"
Set-TextValue "E23" "-1.75%"
Set-TextValue "D24" "0.001220"
Set-TextValue "E24" "-2.62%"
Set-TextValue "D25" "0.004809"
Set-TextValue "E25" "1.77%"
Set-TextValue "D26" "0.0001253"
Set-TextValue "E26" "0.04%"
Set-TextValue "D27" "0.0003008"
Set-TextValue "E27" "-33.52%"
Set-TextValue "D39" "0.01901"
Set-TextValue "E39" "-3.71%"
Set-TextValue "D40" "0.04688"
Set-TextValue "E40" "-4.75%"
Set-TextValue "D41" "0.007536"
Set-TextValue "E41" "-4.51%"
Set-TextValue "D42" "0.009605"
Set-TextValue "E42" "22.58%"
Set-TextValue "D43" "0.1325"
Set-TextValue "E43" "-5.14%"
Set-TextValue "D44" "0.002116"
Set-TextValue "E44" "0.04%"
Set-TextValue "D45" "0.01078"
Set-TextValue "E45" "-8.67%"
Set-TextValue "D46" "0.00006246"
Set-TextValue "E46" "-4.56%"
Set-TextValue "E47" "-0.06%"
Set-TextValue "E48" "84.05%"
Set-TextValue "E49" "-22.03%"
Set-TextValue "E50" "-0.06%"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "-0.06%"
